$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a typo in the existing first table: "M.Sindu" -> "M.Sindhu" ---
$ws.Range("E5").Value = "M.Sindhu"

# --- Add a second "Day 2" table below the first one (rows 10-14) ---
# (Values are written in this particular order so newly-introduced shared
#  strings land in the same append order as the source file.)

$ws.Range("K14").Value = "Simialr products"
$ws.Range("G11").Value = "Wish List dummy data"
$ws.Range("I11").Value = "Not Completed"

# Header row (row 10) - same headers as row 3, bold like the first table's header
$ws.Range("C10").Value = "Days"
$ws.Range("C10").Font.Bold = $true
$ws.Range("E10").Value = "Team Members"
$ws.Range("E10").Font.Bold = $true
$ws.Range("G10").Value = "Task Assigned"
$ws.Range("G10").Font.Bold = $true
$ws.Range("I10").Value = "Status"
$ws.Range("I10").Font.Bold = $true
$ws.Range("K10").Value = "Pending"
$ws.Range("K10").Font.Bold = $true

# Data rows
$ws.Range("C11").Value = 2
$ws.Range("E11").Value = "S.Swaroopa"

$ws.Range("E12").Value = "M.Sindhu"
$ws.Range("G12").Value = "Wish List Page front end"
$ws.Range("I12").Value = "Completed"
$ws.Range("K12").Value = "Linking Part is Pending"

$ws.Range("E13").Value = "K.Ashritha "
$ws.Range("G13").Value = "CapStore Home Page front end"
$ws.Range("I13").Value = "Completed"

$ws.Range("E14").Value = "I.Prakash"
$ws.Range("G14").Value = "Product Page Front End"
$ws.Range("I14").Value = "Partially Completed "

# Give the "Days" number column in the new table the same look as the first
# table (General/center alignment marker used by C4:C7).
$ws.Range("C11").HorizontalAlignment = 1
$ws.Range("C12").HorizontalAlignment = 1
$ws.Range("C13").HorizontalAlignment = 1
$ws.Range("C14").HorizontalAlignment = 1

# --- Update the active cell/selection to match the final editing position ---
$ws.Range("I11").Select() | Out-Null
